# Apply the BOM change: give C8 its own line (10n / C1710) instead of
# lumping it in with the DNP placeholder row, and remove C8 from the
# "100n" group's designator list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: remove "C8" from the list of designators that use the 100n cap.
$ws.Range("B2").Value = "C5, C11, C12, C17"

# Row 4 previously described a DNP (do-not-populate) part shared by C7.
# It now describes C8 specifically, populated with a 10n cap (part C1710).
$ws.Range("A4").Value = "10n"
$ws.Range("B4").Value = "C8"
$ws.Range("C4").Value = "C0805"
$ws.Range("D4").Value = "C1710"
